# Regenerate save_data: use K (strikeouts) instead of Strike# for column G,
# recalculated std/mean, and write updated s_vals into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(1,1,1,2,1,3,0,2,0,2,2,1,0,1,0,3,3,1,4,0,1,1,2,2,1,1,1,1,1,2,1,1,1,0,3,0,1,3,2,0,2,1,1,1,1,1,2,4,2,1,1,0,2,1,2,1,2,3,3,2,2,1,1,0,1,1,3,2)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
